$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: new data row (copied/extended from row 6) -------------------
$ws.Range("A7").Value  = $ws.Range("A6").Value2
$ws.Range("B7").Value  = $ws.Range("B6").Value2
$ws.Range("C7").Value  = $ws.Range("C6").Value2
$ws.Range("D7").Value  = $ws.Range("D6").Value2
$ws.Range("E7").Value  = $ws.Range("E6").Value2
$ws.Range("F7").Value  = $ws.Range("F6").Value2
$ws.Range("G7").Value  = $ws.Range("G6").Value2
$ws.Range("H7").Value  = $ws.Range("H6").Value2
$ws.Range("I7").Value  = $ws.Range("I6").Value2
$ws.Range("J7").Value  = $ws.Range("J6").Value2
$ws.Range("K7").Value  = $ws.Range("K6").Value2

# New unique strings are introduced in this exact order so that the shared
# string table ends up appended in the same order as the authored edit.
$ws.Range("AD7").Value = "malo"

$ws.Range("T2").Value = "xx"
$ws.Range("T3").Value = "yy"
$ws.Range("T4").Value = "zz"
$ws.Range("T5").Value = "ff"
$ws.Range("T6").Value = "gg"
$ws.Range("T7").Value = "dd"

$ws.Range("L2").Value = "quintana.nicolas@javeriana.edi.co"
$ws.Range("L3").Value = "jecheverry@javeriana.edu.co"
$ws.Range("L5").Value = "andrea.torres@gmail.com"
$ws.Range("L6").Value = "valentina.lopez@javeriana.edu.co"

$ws.Range("L7").Value = "mateo.fernandez@.hotmail.com"
$ws.Hyperlinks.Add($ws.Range("L7"), "mailto:mateo.fernandez@.hotmail.com")

# Re-apply the same look-and-feel (underline hyperlink style) that the
# other cells in column L already use, instead of the "new" hyperlink
# style that .Hyperlinks.Add creates by default.
$ws.Range("L6").Copy()
$ws.Range("L7").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("W7").Value = 55555

$ws.Range("AC7").Value = $ws.Range("AC6").Value2
$ws.Range("AM7").Value = $ws.Range("AM6").Value2
$ws.Range("AN7").Value = $ws.Range("AN6").Value2
$ws.Range("AP7").Value = $ws.Range("AP6").Value2
$ws.Range("AQ7").Value = $ws.Range("AQ6").Value2

# --- View state -----------------------------------------------------------
$ws.Activate()
$ws.Range("J22").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
